$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title / Subtitle swap.
#    Before: Title = "Getting Started with R", Subtitle = "R Handout"
#    After : Title = "R Handout",              Subtitle = "Getting Data Into R"
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titleRange = $d.Range($titlePara.Start, $titlePara.End - 1)
$titleRange.Text = "R Handout"

$subtitlePara = $d.Paragraphs(2).Range
$subtitleRange = $d.Range($subtitlePara.Start, $subtitlePara.End - 1)
$subtitleRange.Text = "Getting Data Into R"

# ---------------------------------------------------------------------------
# 2) Rename the example data-frame object from `iris` to `dfobj` throughout
#    the R code blocks (case-sensitive so the literal file name "Iris.csv"
#    used in read.csv(...) is left untouched).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "iris"
$find.Replacement.Text = "dfobj"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $true
$find.MatchWildcards = $false
$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
